$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Observations")

# Insert a new row above row 162 ("rv-age-in-years") for the new
# "rv-age-in-months" CDE, pushing rv-age-in-years (and everything below
# it) down by one row.
$ws.Rows.Item(162).Insert()

# Copy the formatting (style) of the row that now sits at 163 (the
# shifted-down rv-age-in-years row) onto the freshly inserted blank row
# 162, so the new row matches the sheet's standard data-row style.
$ws.Range("A163:K163").Copy()
$ws.Range("A162:K162").PasteSpecial(-4122)

# Populate the new row's values.
$ws.Range("A162").Value2 = "rv-age-in-months"
$ws.Range("B162").Value2 = "RV Age In Months"
$ws.Range("C162").Value2 = "Unified Medical Language System#C0001779, SNOMED CT#397669002, Unified Medical Language System#C2598519, SNOMED CT#424144002, SNOMED CT#125676002, Unified Medical Language System#C0750480, SNOMED CT#246205007"
$ws.Range("E162").Value2 = "null#30525-0"
$ws.Range("F162").Value2 = "http://hl7.org/fhir/ValueSet/observation-codes (example)"
$ws.Range("G162").Value2 = "dateTime, Period, Timing, instant"
$ws.Range("H162").Value2 = "Quantity, CodeableConcept, string, boolean, integer, Range, Ratio, SampledData, time, dateTime, Period"
$ws.Range("I162").Value2 = "optional"
